$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 22 - S03_G02_TB001 (Implement service to convert Alert into Order Request)
$ws.Range("G22").Value = "implemented"
$ws.Range("F22").Value = "Introduced an order creation service that maps Alert fields into an Order in WAITING state."
$ws.Range("H22").Value = "Alert-to-order mapping uses symbol/action/qty/price, defaults to MARKET/MIS and MANUAL mode."
$ws.Range("I22").Value = "Later, incorporate strategy config and risk checks into the transformation."

# Row 23 - S03_G02_TB002 (Persist orders with WAITING status in DB for manual queue)
$ws.Range("G23").Value = "implemented"
$ws.Range("F23").Value = "Webhook now creates a corresponding Order row for each accepted Alert."
$ws.Range("H23").Value = "Orders are stored with status=WAITING and mode=MANUAL, ready for the manual queue."
$ws.Range("I23").Value = "Expose dedicated queue APIs and link to the frontend queue view in Sprint S04."

# Row 24 - S03_G02_TB003 (Expose API to list waiting orders and basic order details for frontend)
$ws.Range("G24").Value = "implemented"
$ws.Range("F24").Value = "Extended webhook tests to assert that Orders are created alongside Alerts."
$ws.Range("H24").Value = "Tests verify the order is linked to the alert and in WAITING/MANUAL state."
$ws.Range("I24").Value = "Add additional tests around failure modes and idempotency if required."
